$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BD7").Value = 126

$ws.Range("G13").Value = 1.75
$ws.Range("H13").Value = 3.1
$ws.Range("I13").Value = 5.75
$ws.Range("J13").Value = 2.5
$ws.Range("L13").Value = 6
$ws.Range("U13").Value = 2.38
$ws.Range("V13").Value = 1.53
$ws.Range("W13").Value = 5
$ws.Range("X13").Value = 7
$ws.Range("Z13").Value = 13
$ws.Range("AD13").Value = 6.5
$ws.Range("AG13").Value = 11
$ws.Range("AH13").Value = 26
$ws.Range("AI13").Value = 19
$ws.Range("AJ13").Value = 67
$ws.Range("AK13").Value = 51
$ws.Range("AL13").Value = 67
$ws.Range("AN13").Value = 3.5
$ws.Range("AO13").Value = 10
$ws.Range("AP13").Value = 26
$ws.Range("AQ13").Value = 34
$ws.Range("AU13").Value = 10
$ws.Range("AW13").Value = 6.5
$ws.Range("AX13").Value = 34
$ws.Range("AZ13").Value = 126
$ws.Range("BA13").Value = 201

$ws.Range("G14").Value = 1.73
$ws.Range("H14").Value = 3.4
$ws.Range("I14").Value = 5.25
$ws.Range("J14").Value = 2.4
$ws.Range("L14").Value = 6
$ws.Range("M14").Value = 1.08
$ws.Range("N14").Value = 8
$ws.Range("Q14").Value = 2.4
$ws.Range("R14").Value = 1.53
$ws.Range("U14").Value = 2.25
$ws.Range("V14").Value = 1.57
$ws.Range("X14").Value = 7
$ws.Range("Z14").Value = 13
$ws.Range("AB14").Value = 41
$ws.Range("AD14").Value = 7
$ws.Range("AG14").Value = 10
$ws.Range("AH14").Value = 23
$ws.Range("AI14").Value = 17
$ws.Range("AN14").Value = 3.5
$ws.Range("AO14").Value = 9.5
$ws.Range("AQ14").Value = 34
$ws.Range("AS14").Value = 251
$ws.Range("AV14").Value = 81
$ws.Range("AW14").Value = 6.5
$ws.Range("AX14").Value = 29
$ws.Range("AZ14").Value = 126

$ws.Range("G15").Value = 2.05
$ws.Range("H15").Value = 3.1
$ws.Range("I15").Value = 3.9
$ws.Range("J15").Value = 2.88
$ws.Range("M15").Value = 1.13
$ws.Range("N15").Value = 6
$ws.Range("O15").Value = 1.53
$ws.Range("P15").Value = 2.38
$ws.Range("V15").Value = 1.57
$ws.Range("AA15").Value = 21
$ws.Range("AN15").Value = 3.75
$ws.Range("AO15").Value = 12
